$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 520, shifting existing rows (and the old last
# row, 587) down by one -- matches the diff's dimension change A1:R587 -> A1:R588.
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row with the new daily record.
$ws.Range("A520").Value = 3
$ws.Range("B520").Value = "Femacal de La Calera"
$ws.Range("C520").Value = "Coquimbo"
$ws.Range("D520").Value = 45124
$ws.Range("E520").Value = 5
$ws.Range("F520").Value = 100114013
$ws.Range("G520").Value = "Zanahoria"
$ws.Range("H520").Value = "Sin especificar"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 230
$ws.Range("K520").Value = 7000
$ws.Range("L520").Value = 7500
$ws.Range("M520").Value = 7348
$ws.Range("N520").Value = "$/saco 20 kilos"
$ws.Range("O520").Value = "Provincia de Quillota"
$ws.Range("P520").Value = 367
$ws.Range("Q520").Value = 20
$ws.Range("R520").Value = "Hortaliza"
